{"js": "// Update the grant-proposal placeholder paragraphs:\n//  - \"Project Name: {{ project_name }}\" -> \"Organization: {{ org_name }}\"\n//  - \"Project Description: {{ project_description }}\" -> \"Project Title: {{ project_title }}\"\n//  - followed by 8 new placeholder paragraphs (Executive Summary .. Budget Overview)\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst nameParagraph = paragraphs.items.find((p) => p.text === \"Project Name: {{ project_name }}\");\nif (nameParagraph) {\n  nameParagraph.insertText(\"Organization: {{ org_name }}\", \"Replace\");\n}\n\nconst descriptionParagraph = paragraphs.items.find(\n  (p) => p.text === \"Project Description: {{ project_description }}\"\n);\nif (descriptionParagraph) {\n  descriptionParagraph.insertText(\"Project Title: {{ project_title }}\", \"Replace\");\n\n  const newParagraphTexts = [\n    \"Executive Summary: {{ summary }}\",\n    \"Background: {{ background }}\",\n    \"Problem: {{ problem }}\",\n    \"Goals & Objectives: {{ goals }}\",\n    \"Project Design: {{ methods }}\",\n    \"Timeline: {{ timeline }}\",\n    \"Evaluation Plan: {{ evaluation }}\",\n    \"Budget Overview: {{ budget }}\",\n  ];\n\n  let anchor = descriptionParagraph;\n  for (const text of newParagraphTexts) {\n    anchor = anchor.insertParagraph(text, \"After\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the grant-proposal placeholder paragraphs:\n#  - \"Project Name: {{ project_name }}\" -> \"Organization: {{ org_name }}\"\n#  - \"Project Description: {{ project_description }}\" -> \"Project Title: {{ project_title }}\"\n#  - followed by 8 new placeholder paragraphs (Executive Summary .. Budget Overview)\n$d = $word.ActiveDocument\n\n$nameRange = $d.Content\n$nameRange.Find.Execute(\n    \"Project Name: {{ project_name }}\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"Organization: {{ org_name }}\",\n    2\n)\n\n$descRange = $d.Content\n$descRange.Find.Execute(\n    \"Project Description: {{ project_description }}\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"Project Title: {{ project_title }}\",\n    2\n)\n\n# Locate the paragraph that now reads \"Project Title: {{ project_title }}\"\n# and append the new placeholder paragraphs right after it, in order.\n$titlePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"Project Title: {{ project_title }}`r\") {\n        $titlePara = $p\n    }\n}\n\n$newTexts = @(\n    \"Executive Summary: {{ summary }}\",\n    \"Background: {{ background }}\",\n    \"Problem: {{ problem }}\",\n    \"Goals & Objectives: {{ goals }}\",\n    \"Project Design: {{ methods }}\",\n    \"Timeline: {{ timeline }}\",\n    \"Evaluation Plan: {{ evaluation }}\",\n    \"Budget Overview: {{ budget }}\"\n)\n\n$cur = $titlePara\nforeach ($t in $newTexts) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $d.Paragraphs($cur.Index + 1)\n    $cur.Range.Text = $t\n}\n"}
